$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "[-, -, Guilherme-Eletrohidráulica, Leandro-Sistemas de Refrigeração]"
$ws.Range("C18").Value = "[Emerson-Eletrônica Básica, Allan Cupertino-Instalções Elétricas]"
$ws.Range("D18").Value = "[Cleidson-Automação Industrial, Paulo Rob.-CAM, Cláudio-Tecnologia da Soldagem, Guilherme-Eletrohidráulica]"
$ws.Range("E18").Value = "Allan Cupertino-Máquinas Elétri"
$ws.Range("F18").Value = "[Allan Cupertino-Instalções Elétricas, Weslei-CAD]"

$ws.Range("B19").Value = "[-, -, Guilherme-Eletrohidráulica, Leandro-Sistemas de Refrigeração]"
$ws.Range("C19").Value = "[Emerson-Eletrônica Básica, Allan Cupertino-Instalções Elétricas]"
$ws.Range("D19").Value = "[Cleidson-Automação Industrial, Paulo Rob.-CAM, Cláudio-Tecnologia da Soldagem, Guilherme-Eletropneumática]"
$ws.Range("E19").Value = "Allan Cupertino-Máquinas Elétri"
$ws.Range("F19").Value = "[Allan Cupertino-Instalções Elétricas, Weslei-CAD]"

$ws.Range("B20").Value = "[-, Guilherme-Eletropneumática, -, Leandro-Sistemas de Refrigeração]"
$ws.Range("C20").Value = "[Allan Cupertino-Lab. De Máquinas elétricas, João Paulo-Lab. Circuitos Elétricos]"
$ws.Range("D20").Value = "[Cleidson-Automação Industrial, Paulo Rob.-CAM, Cláudio-Tecnologia da Soldagem, Guilherme-Eletropneumática]"
$ws.Range("E20").Value = "Andre B.-Circuitos Elétrico"
$ws.Range("F20").Value = "[Emerson-Eletrônica Básica, Weslei-CAD]"

$ws.Range("B21").Value = "[Guilherme-Eletrohidráulica, -, -, Leandro-Sistemas de Refrigeração]"
$ws.Range("C21").Value = "[Allan Cupertino-Lab. De Máquinas elétricas, João Paulo-Lab. Circuitos Elétricos]"
$ws.Range("D21").Value = "[Cleidson-Automação Industrial, Paulo Rob.-CAM, Cláudio-Tecnologia da Soldagem, Guilherme-Eletropneumática]"
$ws.Range("E21").Value = "Andre B.-Circuitos Elétrico"
$ws.Range("F21").Value = "[Emerson-Eletrônica Básica, Weslei-CAD]"
